$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 1742
$ws1.Range("F8").Value = 2166
$ws1.Range("F9").Value = 2083
$ws1.Range("F11").Value = 590
$ws1.Range("F18").Value = 183
$ws1.Range("F19").Value = 1546
$ws1.Range("F20").Value = 594
$ws1.Range("F23").Value = 12131
$ws1.Range("F24").Value = 12142
$ws1.Range("F27").Value = 4
$ws1.Range("F29").Value = 11
$ws1.Range("F30").Value = 324
$ws1.Range("F31").Value = 1903
$ws1.Range("F33").Value = 549

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 73

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 73
$ws4.Range("F6").Value = 1742
$ws4.Range("F9").Value = 2166
$ws4.Range("F10").Value = 2083
$ws4.Range("F12").Value = 590
$ws4.Range("F22").Value = 183
$ws4.Range("F23").Value = 1546
$ws4.Range("F24").Value = 594
$ws4.Range("F27").Value = 12131
$ws4.Range("F28").Value = 12142
$ws4.Range("F31").Value = 4
$ws4.Range("F33").Value = 11
$ws4.Range("F34").Value = 324
$ws4.Range("F35").Value = 1903
$ws4.Range("F39").Value = 549

Write-Output "Done"
